$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.211.75'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '1.859.96'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('D4').Value = '0.9996'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '0.7136'
$ws.Range('E5').Value = '  -0.25%  '
$ws.Range('D6').Value = '240.40'
$ws.Range('D7').Value = '0.9999'
$ws.Range('E7').Value = '  -0.17%  '
$ws.Range('D8').Value = '0.07765'
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('D9').Value = '0.3075'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').Value = '25.13'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').Value = '0.08260'
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = '1.870.90'
$ws.Range('E12').Value = '  -0.27%  '
$ws.Range('D13').Value = '5.230'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '0.7163'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').Value = '90.40'
$ws.Range('D16').Value = '29.204.34'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').Value = '244.66'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('E19').Value = '  -0.73%  '
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').Value = '2.103.05'
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '7.978'
$ws.Range('E23').Value = '  +3.02%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  -0.13%  '
$ws.Range('D25').Value = '0.1593'
$ws.Range('E25').Value = '  +2.80%  '
$ws.Range('D26').Value = '162.44'
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('E27').Value = '  -0.91%  '
$ws.Range('D28').Value = '18.34'
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').Value = '1.495'
$ws.Range('E29').Value = '  +0.92%  '
$ws.Range('E30').Value = '  -3.13%  '
$ws.Range('D31').Value = '4.423'
$ws.Range('E31').Value = '  +2.17%  '
$ws.Range('D32').Value = '4.236'
$ws.Range('E32').Value = '  +3.57%  '
$ws.Range('D33').Value = '0.05185'
$ws.Range('E33').Value = '  -1.11%  '
$ws.Range('D34').Value = '1.912'
$ws.Range('E34').Value = '  -0.75%  '
$ws.Range('E35').Value = '  -2.14%  '
$ws.Range('D36').Value = '0.7272'
$ws.Range('E36').Value = '  +1.55%  '
$ws.Range('D37').Value = '2.673'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').Value = '0.01854'
$ws.Range('E38').Value = '  -0.48%  '
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('D40').Value = '1.165.81'
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('D41').Value = '0.9068'
$ws.Range('E41').Value = '  -0.17%  '
$ws.Range('D42').Value = '6.154'
$ws.Range('E42').Value = '  +2.40%  '
$ws.Range('D43').Value = '72.38'
$ws.Range('E43').Value = '  +0.63%  '
$ws.Range('D44').Value = '0.9994'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').Value = '101.63'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').Value = '2.000.57'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').Value = '0.5218'
$ws.Range('E47').Value = '  -2.83%  '
$ws.Range('E48').Value = '  +0.26%  '
$ws.Range('E49').Value = '  -1.45%  '
$ws.Range('D50').Value = '9.316'
$ws.Range('E50').Value = '  +1.96%  '
$ws.Range('D51').Value = '2.861'
$ws.Range('E51').Value = '  +1.15%  '
